# Updated Script to handle Upload & DownLoad Functionalities
# Record the expected "Success" status for the Tools_QA_Practice test-data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tools_QA_Practice")
$ws.Range("L2").Value = "Success"
